$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names / links) ---
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

# --- Numeric-looking text cells (price / volume%) need Text format so Excel
#     stores the literal string instead of auto-converting to a number/percent ---
$numericLikeCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $numericLikeCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "288.83"
$ws.Range("E2").Value = "-0.05%"
$ws.Range("D3").Value = "30.97"
$ws.Range("E3").Value = "0.87%"
$ws.Range("D4").Value = "4.953"
$ws.Range("E4").Value = "0.31%"
$ws.Range("D5").Value = "0.07373"
$ws.Range("E5").Value = "2.80%"
$ws.Range("D6").Value = "2.325"
$ws.Range("E6").Value = "27.94%"
$ws.Range("D7").Value = "7.719"
$ws.Range("E7").Value = "1.36%"
$ws.Range("D8").Value = "3.719"
$ws.Range("E8").Value = "-0.30%"
$ws.Range("D9").Value = "0.9101"
$ws.Range("E9").Value = "1.52%"
$ws.Range("D10").Value = "0.09111"
$ws.Range("E10").Value = "17.59%"
$ws.Range("D11").Value = "0.1694"
$ws.Range("E11").Value = "0.82%"
$ws.Range("D12").Value = "0.08309"
$ws.Range("E12").Value = "4.60%"
$ws.Range("D13").Value = "0.03120"
$ws.Range("E13").Value = "2.70%"
$ws.Range("D14").Value = "0.09962"
$ws.Range("E14").Value = "-0.46%"
$ws.Range("D15").Value = "0.001492"
$ws.Range("E15").Value = "-0.29%"
$ws.Range("D16").Value = "0.005794"
$ws.Range("E16").Value = "0.21%"
$ws.Range("D17").Value = "3.498"
$ws.Range("E17").Value = "0.96%"
$ws.Range("D18").Value = "2.112"
$ws.Range("E18").Value = "1.81%"
$ws.Range("E19").Value = "0.23%"
$ws.Range("D20").Value = "0.1286"
$ws.Range("E20").Value = "0.30%"
$ws.Range("D21").Value = "4.183"
$ws.Range("E21").Value = "3.62%"
$ws.Range("D22").Value = "0.2097"
$ws.Range("E22").Value = "-0.28%"
$ws.Range("D23").Value = "0.04501"
$ws.Range("E23").Value = "-0.29%"
$ws.Range("D24").Value = "0.001207"
$ws.Range("E24").Value = "-0.57%"
$ws.Range("D25").Value = "0.004174"
$ws.Range("E25").Value = "-9.61%"
$ws.Range("D26").Value = "0.0001297"
$ws.Range("E26").Value = "-0.34%"
$ws.Range("D27").Value = "0.0003389"
$ws.Range("E27").Value = "-95.48%"
$ws.Range("D39").Value = "0.01581"
$ws.Range("E39").Value = "1.18%"
$ws.Range("D40").Value = "0.04464"
$ws.Range("E40").Value = "2.93%"
$ws.Range("D41").Value = "0.007345"
$ws.Range("E41").Value = "-0.16%"
$ws.Range("D42").Value = "0.009496"
$ws.Range("E42").Value = "-5.53%"
$ws.Range("D43").Value = "0.1328"
$ws.Range("E43").Value = "2.16%"
$ws.Range("D44").Value = "0.002256"
$ws.Range("E44").Value = "8.76%"
$ws.Range("D45").Value = "0.008068"
$ws.Range("E45").Value = "-11.71%"
$ws.Range("D46").Value = "0.00006100"
$ws.Range("E46").Value = "2.08%"
$ws.Range("E47").Value = "-0.31%"
$ws.Range("D48").Value = "2.447"
$ws.Range("E48").Value = "8.53%"
$ws.Range("D49").Value = "0.001997"
$ws.Range("E49").Value = "-33.43%"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").Value = "-0.31%"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.31%"
